$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Cells.Item(74, 8).Value = 4815.857  # H74: 5046.2 -> 4815.857
$ws.Cells.Item(74, 9).Value = 4807.75  # I74: 5057.75 -> 4807.75
$ws.Cells.Item(74, 10).Value = 4826.6665  # J74: 5000 -> 4826.6665
$ws.Cells.Item(74, 11).Value = 4807.75  # K74: 5057.75 -> 4807.75
$ws.Cells.Item(74, 12).Value = 4826.6665  # L74: 5000 -> 4826.6665
$ws.Cells.Item(74, 13).Value = -3871.75  # M74: -4121.75 -> -3871.75
$ws.Cells.Item(74, 14).Value = -6698.6665  # N74: -6872 -> -6698.6665

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Cells.Item(76, 8).Value = 3263.3333  # H76: 3336 -> 3263.3333
$ws.Cells.Item(76, 9).Value = 2925  # I76: 2933.3333 -> 2925
$ws.Cells.Item(76, 11).Value = 2925  # K76: 2933.3333 -> 2925
$ws.Cells.Item(76, 13).Value = -2610  # M76: -2618.3333 -> -2610

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Cells.Item(77, 8).Value = 4815.857  # H77: 5046.2 -> 4815.857
$ws.Cells.Item(77, 9).Value = 4807.75  # I77: 5057.75 -> 4807.75
$ws.Cells.Item(77, 10).Value = 4826.6665  # J77: 5000 -> 4826.6665
$ws.Cells.Item(77, 11).Value = 24038.75  # K77: 25288.75 -> 24038.75
$ws.Cells.Item(77, 12).Value = 24133.3325  # L77: 25000 -> 24133.3325
$ws.Cells.Item(77, 13).Value = -19358.75  # M77: -20608.75 -> -19358.75
$ws.Cells.Item(77, 14).Value = -33493.3325  # N77: -34360 -> -33493.3325

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Cells.Item(79, 8).Value = 3263.3333  # H79: 3336 -> 3263.3333
$ws.Cells.Item(79, 9).Value = 2925  # I79: 2933.3333 -> 2925
$ws.Cells.Item(79, 11).Value = 2925  # K79: 2933.3333 -> 2925
$ws.Cells.Item(79, 13).Value = -1833  # M79: -1841.3333 -> -1833

# Row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Cells.Item(106, 8).Value = 1663.6364  # H106: 1570 -> 1663.6364
$ws.Cells.Item(106, 10).Value = 3200  # J106: 3500 -> 3200
$ws.Cells.Item(106, 12).Value = 3200  # L106: 3500 -> 3200
$ws.Cells.Item(106, 14).Value = -4462  # N106: -4762 -> -4462

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Cells.Item(112, 8).Value = 3473.913  # H112: 3022.3333 -> 3473.913
$ws.Cells.Item(112, 10).Value = 4022.2222  # J112: 3326.8 -> 4022.2222
$ws.Cells.Item(112, 12).Value = 12066.6666  # L112: 9980.400000000001 -> 12066.6666
$ws.Cells.Item(112, 14).Value = -14282.6666  # N112: -12196.4 -> -14282.6666

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 3106.6064  # H138: 3047.0862 -> 3106.6064
$ws.Cells.Item(138, 9).Value = 2148.682  # I138: 1875.8889 -> 2148.682
$ws.Cells.Item(138, 10).Value = 3646.9744  # J138: 4067.1614 -> 3646.9744
$ws.Cells.Item(138, 11).Value = 6446.045999999999  # K138: 5627.6667 -> 6446.045999999999
$ws.Cells.Item(138, 12).Value = 10940.9232  # L138: 12201.4842 -> 10940.9232
$ws.Cells.Item(138, 13).Value = -1306.045999999999  # M138: -487.6666999999998 -> -1306.045999999999
$ws.Cells.Item(138, 14).Value = -21220.9232  # N138: -22481.4842 -> -21220.9232

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 23: A Well-rounded Crew | Iron Hoplon
$ws.Cells.Item(23, 8).Value = 0  # H23: 20000 -> 0
$ws.Cells.Item(23, 10).Value = 0  # J23: 20000 -> 0
$ws.Cells.Item(23, 12).Value = 0  # L23: 20000 -> 0
$ws.Cells.Item(23, 14).ClearContents()  # N23: -20518 -> (removed)

# Row 27: Get Me the Hard Stuff | Ironclad Bronze Buckler
$ws.Cells.Item(27, 8).Value = 64005.332  # H27: 37506 -> 64005.332
$ws.Cells.Item(27, 10).Value = 64005.332  # J27: 37506 -> 64005.332
$ws.Cells.Item(27, 12).Value = 64005.332  # L27: 37506 -> 64005.332
$ws.Cells.Item(27, 14).Value = -64373.332  # N27: -37874 -> -64373.332

# Row 32: Ingot We Trust | Steel Ingot
$ws.Cells.Item(32, 8).Value = 8051.19  # H32: 7694.59 -> 8051.19
$ws.Cells.Item(32, 9).Value = 7689.674  # I32: 7132.912 -> 7689.674
$ws.Cells.Item(32, 10).Value = 14920  # J32: 13373.777 -> 14920
$ws.Cells.Item(32, 11).Value = 7689.674  # K32: 7132.912 -> 7689.674
$ws.Cells.Item(32, 12).Value = 14920  # L32: 13373.777 -> 14920
$ws.Cells.Item(32, 13).Value = -7402.674  # M32: -6845.912 -> -7402.674
$ws.Cells.Item(32, 14).Value = -15494  # N32: -13947.777 -> -15494

# Row 96: The Gauntlet Is Cast | High Steel Gauntlets of Fending
$ws.Cells.Item(96, 8).Value = 13364  # H96: 15455 -> 13364
$ws.Cells.Item(96, 10).Value = 13364  # J96: 15455 -> 13364
$ws.Cells.Item(96, 12).Value = 13364  # L96: 15455 -> 13364
$ws.Cells.Item(96, 14).Value = -18856  # N96: -20947 -> -18856

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 1678.8  # H102: 1625 -> 1678.8
$ws.Cells.Item(102, 9).Value = 1235.5  # I102: 1192.7273 -> 1235.5
$ws.Cells.Item(102, 10).Value = 3452  # J102: 3210 -> 3452
$ws.Cells.Item(102, 11).Value = 1235.5  # K102: 1192.7273 -> 1235.5
$ws.Cells.Item(102, 12).Value = 3452  # L102: 3210 -> 3452
$ws.Cells.Item(102, 13).Value = 386.5  # M102: 429.2727 -> 386.5
$ws.Cells.Item(102, 14).Value = -6696  # N102: -6454 -> -6696

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Cells.Item(110, 8).Value = 1629  # H110: 1612.2222 -> 1629
$ws.Cells.Item(110, 9).Value = 1073  # I110: 1000 -> 1073
$ws.Cells.Item(110, 10).Value = 2185  # J110: 1918.3334 -> 2185
$ws.Cells.Item(110, 11).Value = 1073  # K110: 1000 -> 1073
$ws.Cells.Item(110, 12).Value = 2185  # L110: 1918.3334 -> 2185
$ws.Cells.Item(110, 13).Value = 972  # M110: 1045 -> 972
$ws.Cells.Item(110, 14).Value = -6275  # N110: -6008.3334 -> -6275

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Cells.Item(122, 8).Value = 1726.3684  # H122: 2242.5833 -> 1726.3684
$ws.Cells.Item(122, 9).Value = 1566  # I122: 3004 -> 1566
$ws.Cells.Item(122, 10).Value = 1870.7  # J122: 1988.7778 -> 1870.7
$ws.Cells.Item(122, 11).Value = 4698  # K122: 9012 -> 4698
$ws.Cells.Item(122, 12).Value = 5612.1  # L122: 5966.3334 -> 5612.1
$ws.Cells.Item(122, 13).Value = -2248  # M122: -6562 -> -2248
$ws.Cells.Item(122, 14).Value = -10512.1  # N122: -10866.3334 -> -10512.1

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 7018.841  # H132: 7471.878 -> 7018.841
$ws.Cells.Item(132, 9).Value = 7405.773  # I132: 7735.1904 -> 7405.773
$ws.Cells.Item(132, 10).Value = 6631.909  # J132: 7195.4 -> 6631.909
$ws.Cells.Item(132, 11).Value = 22217.319  # K132: 23205.5712 -> 22217.319
$ws.Cells.Item(132, 12).Value = 19895.727  # L132: 21586.2 -> 19895.727
$ws.Cells.Item(132, 13).Value = -19687.319  # M132: -20675.5712 -> -19687.319
$ws.Cells.Item(132, 14).Value = -24955.727  # N132: -26646.2 -> -24955.727

# Row 15: All Ovo That | Iron Skillet
$ws.Cells.Item(15, 8).Value = 5000  # H15: 0 -> 5000
$ws.Cells.Item(15, 10).Value = 5000  # J15: 0 -> 5000
$ws.Cells.Item(15, 12).Value = 5000  # L15: 0 -> 5000
$ws.Cells.Item(15, 14).Value = -5454  # N15: None -> -5454

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal | High Steel Nugget
$ws.Cells.Item(94, 8).Value = 1492.5769  # H94: 1311.4688 -> 1492.5769
$ws.Cells.Item(94, 9).Value = 1634.8334  # I94: 1394.6522 -> 1634.8334
$ws.Cells.Item(94, 10).Value = 1172.5  # J94: 1098.8889 -> 1172.5
$ws.Cells.Item(94, 11).Value = 1634.8334  # K94: 1394.6522 -> 1634.8334
$ws.Cells.Item(94, 12).Value = 1172.5  # L94: 1098.8889 -> 1172.5
$ws.Cells.Item(94, 13).Value = -1183.8334  # M94: -943.6522 -> -1183.8334
$ws.Cells.Item(94, 14).Value = -2074.5  # N94: -2000.8889 -> -2074.5

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Cells.Item(105, 8).Value = 1957.2632  # H105: 1876.7368 -> 1957.2632
$ws.Cells.Item(105, 9).Value = 1794.375  # I105: 1792.5 -> 1794.375
$ws.Cells.Item(105, 10).Value = 2826  # J105: 2326 -> 2826
$ws.Cells.Item(105, 11).Value = 1794.375  # K105: 1792.5 -> 1794.375
$ws.Cells.Item(105, 12).Value = 2826  # L105: 2326 -> 2826
$ws.Cells.Item(105, 13).Value = -47.375  # M105: -45.5 -> -47.375
$ws.Cells.Item(105, 14).Value = -6320  # N105: -5820 -> -6320

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 4979.1665  # H134: 4633.3335 -> 4979.1665
$ws.Cells.Item(134, 9).Value = 4302.125  # I134: 4556 -> 4302.125
$ws.Cells.Item(134, 10).Value = 5520.8  # J134: 4682.5454 -> 5520.8
$ws.Cells.Item(134, 11).Value = 12906.375  # K134: 13668 -> 12906.375
$ws.Cells.Item(134, 12).Value = 16562.4  # L134: 14047.6362 -> 16562.4
$ws.Cells.Item(134, 13).Value = -10371.375  # M134: -11133 -> -10371.375
$ws.Cells.Item(134, 14).Value = -21632.4  # N134: -19117.6362 -> -21632.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall | Elm Lumber
$ws.Cells.Item(22, 8).Value = 301  # H22: 610.2 -> 301
$ws.Cells.Item(22, 9).Value = 301  # I22: 300.5 -> 301
$ws.Cells.Item(22, 10).Value = 0  # J22: 816.6667 -> 0
$ws.Cells.Item(22, 11).Value = 301  # K22: 300.5 -> 301
$ws.Cells.Item(22, 12).Value = 0  # L22: 816.6667 -> 0
$ws.Cells.Item(22, 13).Value = 49  # M22: 49.5 -> 49
$ws.Cells.Item(22, 14).ClearContents()  # N22: -1516.6667 -> (removed)

# Row 99: O Pine | Pine Lumber
$ws.Cells.Item(99, 8).Value = 3564  # H99: 2180.2856 -> 3564
$ws.Cells.Item(99, 9).Value = 3306  # I99: 1552.4 -> 3306
$ws.Cells.Item(99, 10).Value = 3650  # J99: 3750 -> 3650
$ws.Cells.Item(99, 11).Value = 3306  # K99: 1552.4 -> 3306
$ws.Cells.Item(99, 12).Value = 3650  # L99: 3750 -> 3650
$ws.Cells.Item(99, 13).Value = -1808  # M99: -54.40000000000009 -> -1808
$ws.Cells.Item(99, 14).Value = -6646  # N99: -6746 -> -6646

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Cells.Item(122, 8).Value = 71429890  # H122: 62501204 -> 71429890
$ws.Cells.Item(122, 9).Value = 83334210  # I122: 76923930 -> 83334210
$ws.Cells.Item(122, 10).Value = 4000  # J122: 2733.3333 -> 4000
$ws.Cells.Item(122, 11).Value = 250002630  # K122: 230771790 -> 250002630
$ws.Cells.Item(122, 12).Value = 12000  # L122: 8199.999899999999 -> 12000
$ws.Cells.Item(122, 13).Value = -250000180  # M122: -230769340 -> -250000180
$ws.Cells.Item(122, 14).Value = -16900  # N122: -13099.9999 -> -16900

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 3564  # H126: 2180.2856 -> 3564
$ws.Cells.Item(126, 9).Value = 3306  # I126: 1552.4 -> 3306
$ws.Cells.Item(126, 10).Value = 3650  # J126: 3750 -> 3650
$ws.Cells.Item(126, 11).Value = 9918  # K126: 4657.200000000001 -> 9918
$ws.Cells.Item(126, 12).Value = 10950  # L126: 11250 -> 10950
$ws.Cells.Item(126, 13).Value = -7448  # M126: -2187.200000000001 -> -7448
$ws.Cells.Item(126, 14).Value = -15890  # N126: -16190 -> -15890

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33: Cooking with Gas | Chicken Stock
$ws.Cells.Item(33, 8).Value = 140.28572  # H33: 158.66667 -> 140.28572
$ws.Cells.Item(33, 10).Value = 169.9  # J33: 204.875 -> 169.9
$ws.Cells.Item(33, 12).Value = 1019.4  # L33: 1229.25 -> 1019.4
$ws.Cells.Item(33, 14).Value = -1585.4  # N33: -1795.25 -> -1585.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Cells.Item(107, 8).Value = 547.125  # H107: 507.2903 -> 547.125
$ws.Cells.Item(107, 9).Value = 594.35297  # I107: 517.48 -> 594.35297
$ws.Cells.Item(107, 10).Value = 432.42856  # J107: 464.83334 -> 432.42856
$ws.Cells.Item(107, 11).Value = 594.35297  # K107: 517.48 -> 594.35297
$ws.Cells.Item(107, 12).Value = 432.42856  # L107: 464.83334 -> 432.42856
$ws.Cells.Item(107, 13).Value = 1325.64703  # M107: 1402.52 -> 1325.64703
$ws.Cells.Item(107, 14).Value = -4272.42856  # N107: -4304.83334 -> -4272.42856

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 3401.7273  # H113: 2963.2856 -> 3401.7273
$ws.Cells.Item(113, 9).Value = 3313.4443  # I113: 2875 -> 3313.4443
$ws.Cells.Item(113, 10).Value = 3799  # J113: 3081 -> 3799
$ws.Cells.Item(113, 11).Value = 3313.4443  # K113: 2875 -> 3313.4443
$ws.Cells.Item(113, 12).Value = 3799  # L113: 3081 -> 3799
$ws.Cells.Item(113, 13).Value = -1143.4443  # M113: -705 -> -1143.4443
$ws.Cells.Item(113, 14).Value = -8139  # N113: -7421 -> -8139

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Cells.Item(122, 8).Value = 1937.2069  # H122: 2199.6667 -> 1937.2069
$ws.Cells.Item(122, 9).Value = 1639.1  # I122: 1880.875 -> 1639.1
$ws.Cells.Item(122, 10).Value = 2599.6667  # J122: 2837.25 -> 2599.6667
$ws.Cells.Item(122, 11).Value = 4917.299999999999  # K122: 5642.625 -> 4917.299999999999
$ws.Cells.Item(122, 12).Value = 7799.000100000001  # L122: 8511.75 -> 7799.000100000001
$ws.Cells.Item(122, 13).Value = -2467.299999999999  # M122: -3192.625 -> -2467.299999999999
$ws.Cells.Item(122, 14).Value = -12699.0001  # N122: -13411.75 -> -12699.0001

# Row 132: On Board for Lar | Lar Ingot
$ws.Cells.Item(132, 8).Value = 4125  # H132: 3919 -> 4125
$ws.Cells.Item(132, 9).Value = 4945.3335  # I132: 2260.3333 -> 4945.3335
$ws.Cells.Item(132, 10).Value = 3980.2354  # J132: 6407 -> 3980.2354
$ws.Cells.Item(132, 11).Value = 14836.0005  # K132: 6780.999899999999 -> 14836.0005
$ws.Cells.Item(132, 12).Value = 11940.7062  # L132: 19221 -> 11940.7062
$ws.Cells.Item(132, 13).Value = -12306.0005  # M132: -4250.999899999999 -> -12306.0005
$ws.Cells.Item(132, 14).Value = -17000.7062  # N132: -24281 -> -17000.7062

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Cells.Item(7, 8).Value = 3798.5715  # H7: 3982.8572 -> 3798.5715
$ws.Cells.Item(7, 9).Value = 3114.2856  # I7: 3500 -> 3114.2856
$ws.Cells.Item(7, 10).Value = 4482.857  # J7: 4626.6665 -> 4482.857
$ws.Cells.Item(7, 11).Value = 3114.2856  # K7: 3500 -> 3114.2856
$ws.Cells.Item(7, 12).Value = 4482.857  # L7: 4626.6665 -> 4482.857
$ws.Cells.Item(7, 13).Value = -3002.2856  # M7: -3388 -> -3002.2856
$ws.Cells.Item(7, 14).Value = -4706.857  # N7: -4850.6665 -> -4706.857

# Row 46: Supply Side Logic | Boar Leather
$ws.Cells.Item(46, 8).Value = 370946.97  # H46: 625669.4399999999 -> 370946.97
$ws.Cells.Item(46, 9).Value = 468.4375  # I46: 511.1111 -> 468.4375
$ws.Cells.Item(46, 10).Value = 909824.8  # J46: 1429444.4 -> 909824.8
$ws.Cells.Item(46, 11).Value = 468.4375  # K46: 511.1111 -> 468.4375
$ws.Cells.Item(46, 12).Value = 909824.8  # L46: 1429444.4 -> 909824.8
$ws.Cells.Item(46, 13).Value = -280.4375  # M46: -323.1111 -> -280.4375
$ws.Cells.Item(46, 14).Value = -910200.8  # N46: -1429820.4 -> -910200.8

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Cells.Item(93, 8).Value = 1992.4166  # H93: 1901 -> 1992.4166
$ws.Cells.Item(93, 9).Value = 1781.8  # I93: 1552 -> 1781.8
$ws.Cells.Item(93, 10).Value = 2142.8572  # J93: 2250 -> 2142.8572
$ws.Cells.Item(93, 11).Value = 1781.8  # K93: 1552 -> 1781.8
$ws.Cells.Item(93, 12).Value = 2142.8572  # L93: 2250 -> 2142.8572
$ws.Cells.Item(93, 13).Value = -533.8  # M93: -304 -> -533.8
$ws.Cells.Item(93, 14).Value = -4638.8572  # N93: -4746 -> -4638.8572

# Row 122: Hell on Leather | Gaja Leather
$ws.Cells.Item(122, 8).Value = 7235.7915  # H122: 3463.8572 -> 7235.7915
$ws.Cells.Item(122, 9).Value = 10414  # I122: 3465.3333 -> 10414
$ws.Cells.Item(122, 10).Value = 3479.7273  # J122: 3462.75 -> 3479.7273
$ws.Cells.Item(122, 11).Value = 31242  # K122: 10395.9999 -> 31242
$ws.Cells.Item(122, 12).Value = 10439.1819  # L122: 10388.25 -> 10439.1819
$ws.Cells.Item(122, 13).Value = -28792  # M122: -7945.999899999999 -> -28792
$ws.Cells.Item(122, 14).Value = -15339.1819  # N122: -15288.25 -> -15339.1819

# Row 126: Battered Books | Saiga Leather
$ws.Cells.Item(126, 8).Value = 3798.5715  # H126: 3982.8572 -> 3798.5715
$ws.Cells.Item(126, 9).Value = 3114.2856  # I126: 3500 -> 3114.2856
$ws.Cells.Item(126, 10).Value = 4482.857  # J126: 4626.6665 -> 4482.857
$ws.Cells.Item(126, 11).Value = 9342.856800000001  # K126: 10500 -> 9342.856800000001
$ws.Cells.Item(126, 12).Value = 13448.571  # L126: 13879.9995 -> 13448.571
$ws.Cells.Item(126, 13).Value = -6872.856800000001  # M126: -8030 -> -6872.856800000001
$ws.Cells.Item(126, 14).Value = -18388.571  # N126: -18819.9995 -> -18388.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 34: He's Got Legs | Velveteen Sarouel
$ws.Cells.Item(34, 8).Value = 0  # H34: 14000 -> 0
$ws.Cells.Item(34, 10).Value = 0  # J34: 14000 -> 0
$ws.Cells.Item(34, 12).Value = 0  # L34: 14000 -> 0
$ws.Cells.Item(34, 14).ClearContents()  # N34: -14406 -> (removed)
